# Remove some videos from the stimuli list, then add a new one.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rows (1-based) in the original sheet whose values are being removed:
#   9  -> stimuli/Kids_Jack_06.mp4
#   10 -> stimuli/Kids_Jack_07.mp4
#   13 -> stimuli/Kids_Jack_10.mp4
#   17 -> stimuli/Kids_Jack_14.mp4
#   18 -> stimuli/Kids_Jack_15.mp4
#   19 -> stimuli/Kids_Jack_16.mp4
#   30 -> stimuli/Kids_Jack_27.mp4
#   31 -> stimuli/Kids_Jack_28.mp4
#   40 -> stimuli/P08.mp4
#   41 -> stimuli/P08.mpg
# Delete from the bottom up so earlier row numbers stay valid.
$rowsToDelete = @(41, 40, 31, 30, 19, 18, 17, 13, 10, 9)
foreach ($r in $rowsToDelete) {
    $ws.Rows($r).Delete()
}

# Append the new stimulus at the end of the (now shorter) list.
$lastRow = $ws.Cells($ws.Rows.Count, 1).End(-4162).Row
$ws.Cells($lastRow + 1, 1).Value = "stimuli/PairJack.mp4"

# Update the selected cell to match the saved view state.
$ws.Range("C10").Select()
